$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# Original layout:  B1="Unnamed: 0", C1="   SILVER_FOR", D1="SILVER_FOR"
# New layout:       B1="Unnamed: 0.2", C1="Unnamed: 0.1", D1="Unnamed: 0",
#                    E1="   SILVER_FOR", F1="SILVER_FOR"
# D1 already carries the header style (s=1); clone it onto the two new
# trailing header cells (E1, F1) before the old C1/D1 text slides over.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

$ws.Range("F1").Value = "SILVER_FOR"
$ws.Range("E1").Value = "   SILVER_FOR"
$ws.Range("D1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "Unnamed: 0.1"
$ws.Range("B1").Value = "Unnamed: 0.2"

# --- Column A style (index column) ----------------------------------------
# Column A already uses the bordered/bold style (s=1) for rows 2-11; extend
# that same formatting down to the newly added rows 12-21.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A12:A21").PasteSpecial(-4122) | Out-Null

# --- Data grid --------------------------------------------------------------
# Start from a clean slate below the header row: the old D2:D11 values (the
# "SILVER_FOR" series, originally in column D) must not linger once their
# contents move over to column F, and every newly-covered cell in B:F for
# rows 2-21 needs to start blank before being selectively re-populated.
$ws.Range("B2:F21").ClearContents()

# Column A: running index 0..19 for rows 2..21
$a = 0
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value = $a
    $a++
}

# Column B: 0..14 for rows 2..16 (rows 17-21 stay blank)
$b = 0
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 2).Value = $b
    $b++
}

# Column C: 0..9 for rows 2..11 (rows 12-21 stay blank)
$c = 0
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = $c
    $c++
}

# Column D: 0..4 for rows 2..6 (rows 7-21 stay blank)
$d = 0
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = $d
    $d++
}

# Column E ("   SILVER_FOR" series): values for rows 2..6 only
$eVals = @(27.44466, 26.933548, 27.173124, 26.856646, 26.422922)
$r = 2
foreach ($v in $eVals) {
    $ws.Cells.Item($r, 5).Value = $v
    $r++
}

# Column F ("SILVER_FOR" series): values for rows 7..21
$fVals = @(
    30.93059290717292, 30.6353459147918, 30.30023174736436, 29.53532024840348, 29.15916464141611,
    31.76244298992617, 32.23934168175691, 32.5524573182023, 32.38868114596181, 32.10415671664822,
    28.81685345771996, 28.51061683625062, 27.67123246313037, 27.34084839285276, 26.94339908891732
)
$r = 7
foreach ($v in $fVals) {
    $ws.Cells.Item($r, 6).Value = $v
    $r++
}

Write-Output "done"
